$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common values repeated across data rows (row 2 already has these in A,B,C,D,E,H,L)
$colA = "AutoTestAdmin"
$colB = "AutoTestUser"
$colC = "New Transmittal from Automation"
$colD = "UnTick"
$colE = "Correspondence"
$colH = "Test 1 ta.docx"
$colL = "Message for New transmittal"

# Row 2: change F2 to "Issued for Information" and clear M2 (Action-Level2)
$ws.Range("F2").Value = "Issued for Information"
$ws.Range("M2").Value = ""

# Row 3: Issued for Review / Comments for Issued for Review
$ws.Range("A3").Value = $colA
$ws.Range("B3").Value = $colB
$ws.Range("C3").Value = $colC
$ws.Range("D3").Value = $colD
$ws.Range("E3").Value = $colE
$ws.Range("F3").Value = "Issued for Review"
$ws.Range("H3").Value = $colH
$ws.Range("L3").Value = $colL
$ws.Range("M3").Value = "Comments for Issued for Review"

# Row 4: Request for Information / Comments for Request for Information
$ws.Range("A4").Value = $colA
$ws.Range("B4").Value = $colB
$ws.Range("C4").Value = $colC
$ws.Range("D4").Value = $colD
$ws.Range("E4").Value = $colE
$ws.Range("F4").Value = "Request for Information"
$ws.Range("H4").Value = $colH
$ws.Range("L4").Value = $colL
$ws.Range("M4").Value = "Comments for Request for Information"

# Row 5: Issued for Approval / Approved
$ws.Range("A5").Value = $colA
$ws.Range("B5").Value = $colB
$ws.Range("C5").Value = $colC
$ws.Range("D5").Value = $colD
$ws.Range("E5").Value = $colE
$ws.Range("F5").Value = "Issued for Approval"
$ws.Range("H5").Value = $colH
$ws.Range("L5").Value = $colL
$ws.Range("M5").Value = "Approved"
